# Revise section C3 (Space-Time Complexity Analysis) per evaluator's comments.

$d = $word.ActiveDocument

# wdReplaceOne = 1 ; wdReplaceAll = 2

# ---------------------------------------------------------------------------
# Paragraph: "The core algorithm in my program ..." (C3 paragraph 1)
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute(
    "The core algorithm in my program (the PlotDeliveryRoute function) contains a nested for-each loop which operates on the list of undelivered packages. This means the space-time complexity of the routing algorithm can be expressed as O(n^2).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The core algorithm in my program (the PlotDeliveryRoute function) contains a nested for-each loop which operates on the list of undelivered packages. This means the time complexity of the routing algorithm can be expressed as O(n^2). The algorithm sets a list of package addresses on each truck; each list will scale linearly with the number of packages loaded on the truck, giving the algorithm a space complexity of O(n).",
    2) | Out-Null

Write-Host "Step 1 done"

# ---------------------------------------------------------------------------
# Paragraph: "Once a route has been plotted for a truck ..." (C3 paragraph 2)
# This paragraph's first run currently carries a <w:lastRenderedPageBreak/>
# which, per the target, moves earlier in the document (onto the new
# sentence appended to the previous paragraph). Replacing text that spans
# both of the paragraph's original runs merges them and drops the stale
# rendered-page-break marker, matching the target structure.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute(
    "Once a route has been plotted for a truck, a truck will deliver its packages by visiting each delivery address in its " + [char]8220 + "route" + [char]8221 + " list one at a time. This will be a linear operation; or in other words, package delivery can be expressed as O(n).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Once a route has been plotted for a truck, a truck will deliver its packages by visiting each delivery address in its " + [char]8220 + "route" + [char]8221 + " list one at a time. This will be a linear operation; or in other words, the time complexity for package delivery can be expressed as O(n). The space complexity for this section of the program will be the same as that of the core algorithm as they operate on the same lists of package addresses; that is, the space complexity for this portion can be expressed as O(n).",
    2) | Out-Null

Write-Host "Step 2 done"

# ---------------------------------------------------------------------------
# Paragraph: "Packages will be read from a file ..." (C3 paragraph 3)
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute(
    "Packages will be read from a file and initialized one at a time, making it O(n). Similarly, packages will be loaded onto their designated trucks in linear (O(n)) time. As there is a set number of trucks in this scenario, the trucks will be initialized in constant (O(1)) time.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Packages will be read from a file and initialized one at a time, making the space and time complexity of package initialization O(n). Similarly, packages will be loaded onto their designated trucks in linear (O(n)) time. Because every package needs to be loaded onto a truck, the space complexity of this operation will scale linearly as well (O(n) space). As there is a set number of trucks in this scenario, the trucks will be initialized in constant (O(1)) time. The trucks will also take up a constant amount of space, so the space complexity for this operation can be expressed as O(1).",
    2) | Out-Null

Write-Host "Step 3 done"

# ---------------------------------------------------------------------------
# Paragraph: "All operations on the package hash table ..." (C3 paragraph 4)
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute(
    "All operations on the package hash table, including insertion and retrieval, are constant operations and have a time complexity of O(1). The hash table has a space complexity of O(n).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "All operations on the package hash table, including insertion and retrieval, are constant operations and have a time complexity of O(1). The hash table itself has a space complexity of O(n).",
    2) | Out-Null

Write-Host "Step 4 done"

# ---------------------------------------------------------------------------
# Paragraph: "Overall, my program will have a space-time complexity ..." (C3 paragraph 5)
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute(
    "Overall, my program will have a space-time complexity of O(n^2).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Overall, my program will have a time complexity of O(n^2), and a space complexity of O(n).",
    2) | Out-Null

Write-Host "Step 5 done"
